$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Remove the "Classification: Controlled" text-box shapes from the
# first-page and even-page footers.
$fFirst = $sec.Footers.Item(2)
while ($fFirst.Shapes.Count -gt 0) {
    $fFirst.Shapes.Item(1).Delete()
}

$fEven = $sec.Footers.Item(3)
while ($fEven.Shapes.Count -gt 0) {
    $fEven.Shapes.Item(1).Delete()
}
